# Update crypto price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.256.64"
Set-TextValue "E2" "  +0.41%  "
Set-TextValue "D3" "1.589.38"
Set-TextValue "E3" "  +0.73%  "
Set-TextValue "E4" "  -0.17%  "
Set-TextValue "D5" "212.34"
Set-TextValue "E5" "  +1.51%  "
Set-TextValue "D6" "0.501"
Set-TextValue "E7" "  -0.14%  "
Set-TextValue "E8" "  +0.27%  "
Set-TextValue "D9" "0.0608"
Set-TextValue "E9" "  -0.06%  "
Set-TextValue "D10" "19.34"
Set-TextValue "E10" "  -0.81%  "
Set-TextValue "D11" "0.0849"
Set-TextValue "E11" "  +0.54%  "
Set-TextValue "D12" "1.812.77"
Set-TextValue "E12" "  +0.73%  "
Set-TextValue "D13" "1.616.70"
Set-TextValue "E13" "  +2.03%  "
Set-TextValue "D15" "0.520"
Set-TextValue "E15" "  +1.42%  "
Set-TextValue "D16" "64.31"
Set-TextValue "E16" "  -0.18%  "
Set-TextValue "D17" "26.253.64"
Set-TextValue "E17" "  +0.41%  "
Set-TextValue "D18" "0.0₃0726"
Set-TextValue "E18" "  -0.41%  "
Set-TextValue "D19" "7.40"
Set-TextValue "E19" "  +2.25%  "
Set-TextValue "D20" "212.74"
Set-TextValue "E20" "  +2.34%  "
Set-TextValue "E21" "  -0.16%  "
Set-TextValue "E22" "  +1.04%  "
Set-TextValue "D23" "9.01"
Set-TextValue "E23" "  +1.47%  "
Set-TextValue "E24" "  -3.16%  "
Set-TextValue "D25" "144.36"
Set-TextValue "E25" "  +0.30%  "
Set-TextValue "E26" "  -0.17%  "
Set-TextValue "D27" "7.06"
Set-TextValue "E27" "  +1.42%  "
Set-TextValue "E28" "  -0.44%  "
Set-TextValue "D29" "15.18"
Set-TextValue "E29" "  -0.08%  "
Set-TextValue "D30" "0.0498"
Set-TextValue "E30" "  -1.28%  "
Set-TextValue "D31" "1.15"
Set-TextValue "E31" "  +0.91%  "
Set-TextValue "E32" "  -0.33%  "
Set-TextValue "D33" "2.94"
Set-TextValue "E33" "  -0.47%  "
Set-TextValue "D34" "1.334.76"
Set-TextValue "E34" "  +4.46%  "
Set-TextValue "E35" "  -0.95%  "
Set-TextValue "E36" "  -0.51%  "
Set-TextValue "D37" "0.588"
Set-TextValue "E37" "  -3.16%  "
Set-TextValue "E38" "  +0.43%  "
Set-TextValue "D39" "0.818"
Set-TextValue "E39" "  +0.28%  "
Set-TextValue "E40" "  -7.59%  "
Set-TextValue "E41" "  +3.13%  "
Set-TextValue "E42" "  -0.13%  "
Set-TextValue "E43" "  +0.14%  "
Set-TextValue "D44" "0.763"
Set-TextValue "E44" "  -0.17%  "
Set-TextValue "B45" "RocketPoolETH"
Set-TextValue "C45" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D45" "1.726.41"
Set-TextValue "E45" "  +0.76%  "
Set-TextValue "B46" "Aave"
Set-TextValue "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "61.87"
Set-TextValue "E46" "  -0.76%  "
Set-TextValue "D47" "85.67"
Set-TextValue "E47" "  -3.53%  "
Set-TextValue "E48" "  -3.94%  "
Set-TextValue "D49" "0.0502"
Set-TextValue "E49" "  -0.65%  "
Set-TextValue "D50" "0.0973"
Set-TextValue "D51" "1.00"
Set-TextValue "E51" "  -0.44%  "
